# Fruta / hortaliza, semanal
# Update rows 21-22 with new weekly data, and append two more rows (23-24)
# shifting the previously-existing content for row 22 down, per the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21 (new weekly entry) ---
$ws.Range("A21").Value = 6
$ws.Range("B21").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C21").Value = "Metropolitana"
$ws.Range("D21").Value = 44748
$ws.Range("E21").Value = 13
$ws.Range("F21").Value = 100112035
$ws.Range("G21").Value = "Bruselas (repollito)"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 200
$ws.Range("K21").Value = 16000
$ws.Range("L21").Value = 17000
$ws.Range("M21").Value = 16400
$ws.Range("N21").Value = "$/malla 15 kilos"
$ws.Range("O21").Value = "Provincia de Quillota"
$ws.Range("P21").Value = 1093
$ws.Range("Q21").Value = 15
$ws.Range("R21").Value = "Hortaliza"

# --- Row 22 (previously row 21's content) ---
$ws.Range("A22").Value = 6
$ws.Range("B22").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C22").Value = "Metropolitana"
$ws.Range("D22").Value = 44699
$ws.Range("E22").Value = 13
$ws.Range("F22").Value = 100112035
$ws.Range("G22").Value = "Bruselas (repollito)"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 150
$ws.Range("K22").Value = 18000
$ws.Range("L22").Value = 20000
$ws.Range("M22").Value = 18667
$ws.Range("N22").Value = "$/malla 15 kilos"
$ws.Range("O22").Value = "Provincia de Quillota"
$ws.Range("P22").Value = 1244
$ws.Range("Q22").Value = 15
$ws.Range("R22").Value = "Hortaliza"

# --- Row 23 (new weekly entry) ---
$ws.Range("A23").Value = 6
$ws.Range("B23").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C23").Value = "Metropolitana"
$ws.Range("D23").Value = 44747
$ws.Range("E23").Value = 13
$ws.Range("F23").Value = 100112035
$ws.Range("G23").Value = "Bruselas (repollito)"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 400
$ws.Range("K23").Value = 17000
$ws.Range("L23").Value = 19000
$ws.Range("M23").Value = 17850
$ws.Range("N23").Value = "$/malla 15 kilos"
$ws.Range("O23").Value = "Provincia de Quillota"
$ws.Range("P23").Value = 1190
$ws.Range("Q23").Value = 15
$ws.Range("R23").Value = "Hortaliza"

# --- Row 24 (previously row 22's content) ---
$ws.Range("A24").Value = 6
$ws.Range("B24").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C24").Value = "Metropolitana"
$ws.Range("D24").Value = 44483
$ws.Range("E24").Value = 13
$ws.Range("F24").Value = 100112035
$ws.Range("G24").Value = "Bruselas (repollito)"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 220
$ws.Range("K24").Value = 18000
$ws.Range("L24").Value = 20000
$ws.Range("M24").Value = 18909
$ws.Range("N24").Value = "$/malla 15 kilos"
$ws.Range("O24").Value = "Provincia de Quillota"
$ws.Range("P24").Value = 1261
$ws.Range("Q24").Value = 15
$ws.Range("R24").Value = "Hortaliza"

# Match style of date column (D) for the two new rows, same as the existing date cells.
$ws.Range("D23").NumberFormat = $ws.Range("D22").NumberFormat
$ws.Range("D24").NumberFormat = $ws.Range("D22").NumberFormat
